# Apply the "update models and document" edit to slide 9 (SDK Tool Flow):
#  - "IEC 61499 Models"  -> "IEC 61499 " + "FBTypes" (two runs)
#  - "OPC-UA Models"     -> paragraph 1: "OPC-UA " + "Info" (two runs)
#                           paragraph 2 (new): "Models"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# --- Shape "Snip Single Corner Rectangle 9" (IEC 61499 Models) ---
$iecShape = $s.Shapes.Item(6)
$iecRange = $iecShape.TextFrame.TextRange
$iecRange.Text = "IEC 61499 "
$iecRange.InsertAfter("FBTypes")

# --- Shape "Snip Single Corner Rectangle 10" (OPC-UA Models) ---
$opcShape = $s.Shapes.Item(7)
$opcRange = $opcShape.TextFrame.TextRange
$opcRange.Text = "OPC-UA "
$opcRange.InsertAfter("Info")
$opcRange2 = $opcShape.TextFrame.TextRange
$opcRange2.InsertAfter([char]13 + "Models")
